$d = $word.ActiveDocument

# 1. "sin tener" was split across two runs ("ten" + "er") separated by a
#    _GoBack bookmark; Word recombines it into a single run "tener" once
#    the text is re-typed/edited as a whole.
$d.Content.Find.Execute("eventos sin tener que", $false, $false, $false, $false, $false, $true, 1, $false, "eventos sin tener que", 2) | Out-Null

# 2. "En un futuro" -> "En el futuro"
$d.Content.Find.Execute("En un futuro", $false, $false, $false, $false, $false, $true, 1, $false, "En el futuro", 2) | Out-Null
